$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Clear the content (but keep formatting) of rows that were removed in the edit:
# Row 36: "Supervision Infra Info PTS" entry removed
$ws.Range("D36").Hyperlinks.Delete()
$ws.Range("A36").ClearContents()
$ws.Range("C36:D36").ClearContents()

# Row 63: "WMOS 2015 ESL" entry removed
$ws.Range("D63").Hyperlinks.Delete()
$ws.Range("A63").ClearContents()
$ws.Range("C63:D63").ClearContents()

# Row 70: "WMOS 2015 ESL" entry removed
$ws.Range("D70").Hyperlinks.Delete()
$ws.Range("A70").ClearContents()
$ws.Range("C70:D70").ClearContents()

# Add "ignore" marker in column E for the "IMS" rows that didn't have it yet
$ws.Range("E43").Value = "ignore"
$ws.Range("E50").Value = "ignore"
$ws.Range("E62").Value = "ignore"
$ws.Range("E69").Value = "ignore"
$ws.Range("E74").Value = "ignore"
$ws.Range("E79").Value = "ignore"

$ws.Range("A7").Select()
$ws.Range("A36").Select()
